$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet from "example1" to "Sheet1"
$ws.Name = "Sheet1"

# Clean up header labels: drop the dot-prefixed duplicate-name suffixes added by
# R's make.names()/make.unique() (e.g. "sample.one" -> "sample one"), and drop the
# leading "X" that gets prepended to purely-numeric column names ("X9" -> "9").
$headers = @(
    "sample one",
    "sample one.1",
    "sample two",
    "NA",
    "sample three",
    "sample three.1",
    "9",
    "10",
    "11",
    "12"
)

$headerRange = $ws.Range("A1:J1")

# Force text storage so the purely-numeric headers ("9".."12") stay strings
# instead of being auto-converted to numbers.
$headerRange.NumberFormat = "@"

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Make the header row stand out: bold + centered.
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108  # xlCenter
